$d = $word.ActiveDocument

# The "Requisitos" bullet list currently lists:
#   1) LOB1056 - ... (Requisito fraco)
#   2) LOQ4002 - ... (Requisito fraco)
#   3) LOQ4057 - ... (Requisito fraco)
# It needs to become:
#   1) LOQ4057 - ... (Requisito fraco)
#   2) LOB1056 - ... (Requisito fraco)
#   3) LOQ4002 - ... (Requisito fraco)
# i.e. move the LOQ4057 line (including its trailing line break) so that it
# appears right before the LOB1056 line.

$loq4057Text = "LOQ4057 -  Operações Unitárias III  (Requisito fraco)"
$lob1056Text = "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)"

# Locate the LOQ4057 line (text run) in the document.
$found = $d.Content
$found.Find.Execute($loq4057Text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Extend the found range by one character so it also captures the manual
# line break (w:br) that terminates the run.
$fullLine = $d.Range($found.Start, $found.End + 1)
$lineText = $fullLine.Text

# Remove the LOQ4057 line (with its trailing break) from its current
# location at the end of the list.
$fullLine.Delete()

# Find where the LOB1056 line now starts and insert the LOQ4057 line right
# before it, restoring the trailing line break character that Delete/Insert
# preserve as part of the copied text.
$target = $d.Content
$target.Find.Execute($lob1056Text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertionPoint = $d.Range($target.Start, $target.Start)
$insertionPoint.InsertBefore($lineText)
